# Auto-generated update of cryptos sheet values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.09"
$ws.Range("E2").Value = "'-3.63%"
$ws.Range("D3").Value = "'49.65"
$ws.Range("E3").Value = "'1.12%"
$ws.Range("D4").Value = "'5.164"
$ws.Range("E4").Value = "'-1.75%"
$ws.Range("D5").Value = "'0.07772"
$ws.Range("E5").Value = "'-3.92%"
$ws.Range("D6").Value = "'4.516"
$ws.Range("E6").Value = "'-2.26%"
$ws.Range("D7").Value = "'1.374"
$ws.Range("E7").Value = "'13.95%"
$ws.Range("E8").Value = "'-5.98%"
$ws.Range("E9").Value = "'-6.22%"
$ws.Range("D10").Value = "'0.2003"
$ws.Range("E10").Value = "'2.87%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09557"
$ws.Range("E11").Value = "'0.42%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04727"
$ws.Range("E12").Value = "'5.84%"
$ws.Range("E13").Value = "'-0.13%"
$ws.Range("D14").Value = "'0.001266"
$ws.Range("E14").Value = "'-4.65%"
$ws.Range("D15").Value = "'0.04170"
$ws.Range("E15").Value = "'-3.31%"
$ws.Range("D16").Value = "'0.005836"
$ws.Range("E16").Value = "'-1.17%"
$ws.Range("E17").Value = "'2,017.58%"
$ws.Range("D18").Value = "'3.336"
$ws.Range("E18").Value = "'-0.83%"
$ws.Range("D19").Value = "'2.238"
$ws.Range("E19").Value = "'-8.05%"
$ws.Range("D20").Value = "'0.3486"
$ws.Range("E20").Value = "'2.79%"
$ws.Range("D21").Value = "'7.923"
$ws.Range("E21").Value = "'-3.86%"
$ws.Range("D22").Value = "'0.1338"
$ws.Range("E22").Value = "'-5.26%"
$ws.Range("D23").Value = "'0.3033"
$ws.Range("E23").Value = "'3.80%"
$ws.Range("D24").Value = "'0.001274"
$ws.Range("E24").Value = "'-2.69%"
$ws.Range("D25").Value = "'0.004048"
$ws.Range("E25").Value = "'-4.59%"
$ws.Range("E26").Value = "'-0.24%"
$ws.Range("D38").Value = "'0.02605"
$ws.Range("E38").Value = "'-3.84%"
$ws.Range("D39").Value = "'0.05857"
$ws.Range("E39").Value = "'5.06%"
$ws.Range("E40").Value = "'69.91%"
$ws.Range("D41").Value = "'0.007950"
$ws.Range("E41").Value = "'3.58%"
$ws.Range("D42").Value = "'0.1423"
$ws.Range("E42").Value = "'-1.04%"
$ws.Range("D43").Value = "'0.008424"
$ws.Range("E43").Value = "'9.20%"
$ws.Range("D44").Value = "'0.007645"
$ws.Range("E44").Value = "'-5.64%"
$ws.Range("D45").Value = "'0.3401"
$ws.Range("E45").Value = "'6.50%"
$ws.Range("D46").Value = "'0.00007020"
$ws.Range("E46").Value = "'0.38%"
$ws.Range("E47").Value = "'-0.19%"
$ws.Range("D48").Value = "'0.05230"
$ws.Range("E48").Value = "'-14.72%"
$ws.Range("D49").Value = "'0.002619"
$ws.Range("E49").Value = "'-34.64%"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("E51").Value = "'-0.19%"
